$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fix characteristic names: NH3 -> Ammonia, NO3 -> Nitrate
# (Ammonia is entered first so it lands earlier in the shared strings table)
$ws.Range("A11").Value = "Ammonia"
$ws.Range("A10").Value = "Nitrate"

# Leave the final selection on A11, matching the author's last edit
$ws.Range("A11").Select()
